$d = $word.ActiveDocument

# Replace all occurrences of "July 01, 2022" with "July 02, 2022" (covers the
# three separate spots: change-of-plea date, fine-paid-in-full date, and the
# "license is suspended from" date).
$d.Content.Find.Execute("July 01, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "July 02, 2022", 2)

# Replace the community-control completion deadline date.
$d.Content.Find.Execute("August 30, 2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "August 31, 2022", 2)
